$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2024")

# Insert a new row at row 36, shifting existing rows 36:128 down to 37:129
$ws.Rows.Item(36).Insert()

# Populate the new row's September columns (R = Details, S = Date)
$ws.Range("R36").Value = "reward points cash"
$ws.Range("S36").Value = "2024-09-10 19:43:35"
